$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text that Excel would otherwise auto-convert to a
# number (dropping trailing zeros / switching to scientific notation), so
# we write it with a leading apostrophe to force text, then restore the
# cell's original style (the apostrophe flags the style as quote-prefixed).

$style = $ws.Range("D2").Style
$ws.Range("D2").Value = "'34.838.87"
$ws.Range("D2").Style = $style
$ws.Range("E2").Value = "  -1.56%  "
$style = $ws.Range("D3").Style
$ws.Range("D3").Value = "'1.876.52"
$ws.Range("D3").Style = $style
$ws.Range("E3").Value = "  -2.32%  "
$ws.Range("E4").Value = "  -0.76%  "
$style = $ws.Range("D5").Style
$ws.Range("D5").Value = "'247.06"
$ws.Range("D5").Style = $style
$ws.Range("E5").Value = "  -2.41%  "
$style = $ws.Range("D6").Style
$ws.Range("D6").Value = "'0.685"
$ws.Range("D6").Style = $style
$ws.Range("E6").Value = "  -7.17%  "
$ws.Range("E7").Value = "  -0.80%  "
$style = $ws.Range("D8").Style
$ws.Range("D8").Value = "'42.00"
$ws.Range("D8").Style = $style
$ws.Range("E8").Value = "  +3.16%  "
$style = $ws.Range("D9").Style
$ws.Range("D9").Value = "'0.347"
$ws.Range("D9").Style = $style
$ws.Range("E9").Value = "  -2.75%  "
$style = $ws.Range("D10").Style
$ws.Range("D10").Value = "'51.05"
$ws.Range("D10").Style = $style
$ws.Range("E10").Value = "  -2.94%  "
$style = $ws.Range("D11").Style
$ws.Range("D11").Value = "'0.0738"
$ws.Range("D11").Style = $style
$ws.Range("E11").Value = "  +0.17%  "
$style = $ws.Range("D12").Style
$ws.Range("D12").Value = "'0.0970"
$ws.Range("D12").Style = $style
$ws.Range("E12").Value = "  -2.83%  "
$style = $ws.Range("D13").Style
$ws.Range("D13").Value = "'12.91"
$ws.Range("D13").Style = $style
$ws.Range("E13").Value = "  +1.98%  "
$style = $ws.Range("D14").Style
$ws.Range("D14").Value = "'2.149.19"
$ws.Range("D14").Style = $style
$ws.Range("E14").Value = "  -2.27%  "
$style = $ws.Range("D15").Style
$ws.Range("D15").Value = "'0.717"
$ws.Range("D15").Style = $style
$ws.Range("E15").Value = "  +0.18%  "
$style = $ws.Range("D16").Style
$ws.Range("D16").Value = "'4.90"
$ws.Range("D16").Style = $style
$ws.Range("E16").Value = "  +0.05%  "
$style = $ws.Range("D17").Style
$ws.Range("D17").Value = "'1.874.60"
$ws.Range("D17").Style = $style
$ws.Range("E17").Value = "  -2.62%  "
$style = $ws.Range("D18").Style
$ws.Range("D18").Value = "'34.819.52"
$ws.Range("D18").Style = $style
$ws.Range("E18").Value = "  -1.63%  "
$style = $ws.Range("D19").Style
$ws.Range("D19").Value = "'72.84"
$ws.Range("D19").Style = $style
$style = $ws.Range("D20").Style
$ws.Range("D20").Value = "'0.0₃0821"
$ws.Range("D20").Style = $style
$ws.Range("E20").Value = "  -1.51%  "
$style = $ws.Range("D21").Style
$ws.Range("D21").Value = "'244.13"
$ws.Range("D21").Style = $style
$ws.Range("E21").Value = "  +0.98%  "
$style = $ws.Range("D22").Style
$ws.Range("D22").Value = "'12.74"
$ws.Range("D22").Style = $style
$ws.Range("E22").Value = "  -2.32%  "
$style = $ws.Range("D23").Style
$ws.Range("D23").Value = "'4.93"
$ws.Range("D23").Style = $style
$ws.Range("E23").Value = "  -2.93%  "
$ws.Range("E24").Value = "  -0.75%  "
$style = $ws.Range("D25").Style
$ws.Range("D25").Value = "'2.46"
$ws.Range("D25").Style = $style
$ws.Range("E25").Value = "  +3.78%  "
$ws.Range("E26").Value = "  -3.40%  "
$style = $ws.Range("D27").Style
$ws.Range("D27").Value = "'165.23"
$ws.Range("D27").Style = $style
$ws.Range("E27").Value = "  -1.51%  "
$style = $ws.Range("D28").Style
$ws.Range("D28").Value = "'8.39"
$ws.Range("D28").Style = $style
$ws.Range("E28").Value = "  -3.81%  "
$style = $ws.Range("D29").Style
$ws.Range("D29").Value = "'18.22"
$ws.Range("D29").Style = $style
$ws.Range("E29").Value = "  -3.21%  "
$ws.Range("E30").Value = "  -6.24%  "
$style = $ws.Range("D31").Style
$ws.Range("D31").Value = "'4.128.38"
$ws.Range("D31").Style = $style
$ws.Range("E31").Value = "  -0.08%  "
$style = $ws.Range("D32").Style
$ws.Range("D32").Value = "'1.69"
$ws.Range("D32").Style = $style
$ws.Range("E32").Value = "  +0.47%  "
$ws.Range("E33").Value = "  -1.68%  "
$style = $ws.Range("D34").Style
$ws.Range("D34").Value = "'0.0578"
$ws.Range("D34").Style = $style
$ws.Range("E34").Value = "  -0.63%  "
$style = $ws.Range("D35").Style
$ws.Range("D35").Value = "'4.16"
$ws.Range("D35").Style = $style
$ws.Range("E35").Value = "  -2.23%  "
$ws.Range("E36").Value = "  -0.82%  "
$style = $ws.Range("D37").Style
$ws.Range("D37").Value = "'0.827"
$ws.Range("D37").Style = $style
$ws.Range("E37").Value = "  -8.96%  "
$style = $ws.Range("D38").Style
$ws.Range("D38").Value = "'1.99"
$ws.Range("D38").Style = $style
$ws.Range("E38").Value = "  -1.86%  "
$style = $ws.Range("D40").Style
$ws.Range("D40").Value = "'97.95"
$ws.Range("D40").Style = $style
$ws.Range("E40").Value = "  -1.10%  "
$style = $ws.Range("D41").Style
$ws.Range("D41").Value = "'16.90"
$ws.Range("D41").Style = $style
$ws.Range("E42").Value = "  +0.76%  "
$ws.Range("E43").Value = "  +0.24%  "
$ws.Range("E44").Value = "  -4.95%  "
$style = $ws.Range("D45").Style
$ws.Range("D45").Value = "'1.284.01"
$ws.Range("D45").Style = $style
$ws.Range("E45").Value = "  -4.72%  "
$style = $ws.Range("D46").Style
$ws.Range("D46").Value = "'2.33"
$ws.Range("D46").Style = $style
$ws.Range("E46").Value = "  -6.24%  "
$style = $ws.Range("D47").Style
$ws.Range("D47").Value = "'0.0799"
$ws.Range("D47").Style = $style
$ws.Range("E47").Value = "  +9.40%  "
$ws.Range("E48").Value = "  -1.02%  "
$ws.Range("E49").Value = "  -1.91%  "
$style = $ws.Range("D50").Style
$ws.Range("D50").Value = "'12.18"
$ws.Range("D50").Style = $style
$ws.Range("E50").Value = "  +5.80%  "
$ws.Range("E51").Value = "  -4.03%  "
